# Auto-generated Excel COM-interop edit script
# Applies cell value changes to Sheet1 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("F3").Value = 2.46
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 2.36
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.03
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.01

# Row 4
$ws.Range("F4").Value = 3.65
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 1.45
$ws.Range("I4").Value = 2.72
$ws.Range("J4").Value = 2.8
$ws.Range("K4").Value = 3.25
$ws.Range("N4").Value = 1.54
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 1.25
$ws.Range("Q4").Value = 1.55
$ws.Range("T4").Value = 1.01
$ws.Range("U4").Value = 1.01
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 1000
$ws.Range("AK4").Value = 1000
$ws.Range("AL4").Value = 1000
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 1000

# Row 5
$ws.Range("F5").Value = 1.04
$ws.Range("G5").Value = 980
$ws.Range("H5").Value = 1.04
$ws.Range("I5").Value = 980
$ws.Range("J5").Value = 1.01
$ws.Range("K5").Value = 980

# Row 6
$ws.Range("F6").Value = 1.04
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 1.04
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1.01
$ws.Range("K6").Value = 980
$ws.Range("P6").Value = 1.24
$ws.Range("Q6").Value = 1.01

# Row 8
$ws.Range("F8").Value = 1.04
$ws.Range("G8").Value = 1000
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.01
$ws.Range("K8").Value = 1000
$ws.Range("P8").Value = 1.24
$ws.Range("Q8").Value = 1.01

# Row 9
$ws.Range("F9").Value = 1.35
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 2.04
$ws.Range("I9").Value = 2.28
$ws.Range("J9").Value = 1.33
$ws.Range("P9").Value = 1.24

# Row 10
$ws.Range("H10").Value = 2.26
$ws.Range("J10").Value = 3.6
$ws.Range("K10").Value = 3.7
$ws.Range("M10").Value = 1.06
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 1.06
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.01
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AB10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AD10").Value = 1000
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000

# Row 11
$ws.Range("F11").Value = 4.2
$ws.Range("G11").Value = 980
$ws.Range("H11").Value = 1.87
$ws.Range("J11").Value = 3.25
$ws.Range("P11").Value = 1.25
$ws.Range("Q11").Value = 1.07
$ws.Range("T11").Value = 1.01
$ws.Range("U11").Value = 1.01
$ws.Range("Y11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000

# Row 12
$ws.Range("G12").Value = 1.36
$ws.Range("H12").Value = 1.05
$ws.Range("I12").Value = 30
$ws.Range("J12").Value = 1.05
$ws.Range("N12").Value = 1.03
$ws.Range("P12").Value = 1.25
$ws.Range("Q12").Value = 1.22
$ws.Range("R12").Value = 1.24
$ws.Range("S12").Value = 1.01
$ws.Range("T12").Value = 1.01
$ws.Range("U12").Value = 1.01
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 1000
$ws.Range("AA12").Value = 1000
$ws.Range("AB12").Value = 1000
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("AH12").Value = 1000
$ws.Range("AI12").Value = 1000
$ws.Range("AJ12").Value = 1000
$ws.Range("AL12").Value = 1000
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000

# Row 13
$ws.Range("I13").Value = 980
$ws.Range("J13").Value = 3.95
$ws.Range("P13").Value = 1.25
$ws.Range("Q13").Value = 1.06
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("X13").Value = 1000
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 1000
$ws.Range("AC13").Value = 1000
$ws.Range("AD13").Value = 1000
$ws.Range("AF13").Value = 1000
$ws.Range("AG13").Value = 9.8
$ws.Range("AI13").Value = 1000
$ws.Range("AJ13").Value = 1000
$ws.Range("AK13").Value = 1000
$ws.Range("AM13").Value = 1000

# Row 15
$ws.Range("F15").Value = 4.6
$ws.Range("G15").Value = 1000
$ws.Range("J15").Value = 4.2
$ws.Range("K15").Value = 8.8
$ws.Range("N15").Value = 1.03
$ws.Range("P15").Value = 1.25
$ws.Range("Q15").Value = 1.23
$ws.Range("R15").Value = 1.24
$ws.Range("S15").Value = 1.01
$ws.Range("T15").Value = 1.01
$ws.Range("U15").Value = 1.01
$ws.Range("Y15").Value = 1000
$ws.Range("Z15").Value = 1000
$ws.Range("AB15").Value = 1000
$ws.Range("AC15").Value = 1000
$ws.Range("AD15").Value = 1000
$ws.Range("AE15").Value = 1000
$ws.Range("AG15").Value = 1000
$ws.Range("AJ15").Value = 160
$ws.Range("AO15").Value = 1000

# Row 16
$ws.Range("G16").Value = 2.4
$ws.Range("I16").Value = 4.1
$ws.Range("J16").Value = 3.2
$ws.Range("K16").Value = 3.35
$ws.Range("P16").Value = 1.25

# Row 17
$ws.Range("H17").Value = 7.4
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 60
$ws.Range("P17").Value = 1.25

# Row 19
$ws.Range("F19").Value = 2.9
$ws.Range("G19").Value = 3.3
$ws.Range("H19").Value = 2.56
$ws.Range("I19").Value = 2.88
$ws.Range("J19").Value = 3.2
$ws.Range("P19").Value = 1.25
$ws.Range("Q19").Value = 1.01

# Row 20
$ws.Range("G20").Value = 2.22
$ws.Range("H20").Value = 3.7
$ws.Range("I20").Value = 4.6
$ws.Range("K20").Value = 3.6
$ws.Range("P20").Value = 1.07
$ws.Range("Q20").Value = 1.01
